$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.904.57"
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").Value = "1.709.54"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9976"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3750"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.37"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3448"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.218"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07544"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9987"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.315"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.080"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.81%  "
$ws.Range("D16").Value = "1.705.91"
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001132"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9978"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "84.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.395"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.55%  "
$ws.Range("D24").Value = "24.920.19"
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.442"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.798"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "133.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.86%  "
$ws.Range("D30").Value = "1.894.31"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.235"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +26.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.869"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.234"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.44%  "
$ws.Range("B35").Value = "WEMIXTOKEN"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.778"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.34%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08799"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.648"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06662"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.191"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02414"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2233"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.57%  "
$ws.Range("E42").Value = "  +1.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6461"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9978"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  +4.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6160"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.822"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.129"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07321"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.43%  "
